# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
#
# The "Price" column (D) stores plain text values such as "1.003" or
# "1.290" (note the significant trailing zero) rather than numbers, so
# any numeric-looking replacement is written with a leading apostrophe
# to force Excel to keep it as text instead of silently re-parsing it
# into a Number cell (which would also eat meaningful trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.934.42"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.812.12"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'310.73"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4975"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("D8").Value = "'0.3933"
$ws.Range("D9").Value = "'0.09606"
$ws.Range("E9").Value = "  +23.65%  "
$ws.Range("D10").Value = "'1.101"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "'40.92"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'6.433"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "'1.003"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'20.44"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "1.810.86"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "'7.279"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "'0.00001123"
$ws.Range("E17").Value = "  +4.91%  "
$ws.Range("D18").Value = "'92.36"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'0.06655"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'17.12"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'5.911"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "27.990.85"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'11.14"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "'2.253"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "'159.56"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "2.021.57"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'2.384"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").Value = "'128.22"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").Value = "'0.1066"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "'5.565"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Value = "'0.06704"
$ws.Range("E35").Value = "  -5.19%  "
$ws.Range("D36").Value = "'8.937"
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("D37").Value = "'0.02325"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'0.2137"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "'4.935"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'11.22"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").Value = "'0.6168"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D43").Value = "'1.146"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'13.19"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.290"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5879"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "'3.693"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "'123.05"
$ws.Range("D49").Value = "'1.935"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "'1.177"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "'0.06771"
$ws.Range("E51").Value = "  +0.04%  "
